$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header E1 from "remark" to "Expected "
$ws.Range("E1").Value = "Expected "

# Apply bold font + yellow fill to header row A1:E1, using a template cell +
# paste-formats so only a single new style is minted (matches a real
# "format painter" workflow and avoids extra intermediate styles).
$tmpl = $ws.Range("Z1")
$tmpl.Font.Bold = $true
$tmpl.Interior.Color = 65535
$tmpl.Copy() | Out-Null
$ws.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$tmpl.Clear() | Out-Null
$excel.CutCopyMode = 0

# Column widths (COM ColumnWidth quantizes to whole pixels for the sheet's
# default font, i.e. steps of 1/6 character; subtract the ~5/6-character
# padding offset that Excel's pixel<->XML-width conversion applies so the
# post-quantization stored width lands as close as possible to the target).
$ws.Columns.Item(1).ColumnWidth = 16.619791666666668
$ws.Columns.Item(3).ColumnWidth = 11.529947916666666
$ws.Columns.Item(4).ColumnWidth = 13.893229166666666

# Selection
$ws.Range("B11").Select()

# Page setup - portrait orientation
$ws.PageSetup.Orientation = 1
